# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45172 to serial date 45175.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("C2:C423")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45172) {
        $cell.Value = 45175
    }
}
